$wb = $excel.ActiveWorkbook

# --- Sheet "Alt2": collapse selection from C16:C24 to just C16 ---
$wsAlt2 = $wb.Worksheets.Item("Alt2")
$wsAlt2.Activate()
$wsAlt2.Range("C16").Select()

# --- Sheet "EDA": add sign-handling helper columns (G:M) ---
$ws = $wb.Worksheets.Item("EDA")
$ws.Activate()

# spacer columns G and L get a narrow width + light-gray fill style
$ws.Columns("G").ColumnWidth = 2.5546875
$ws.Columns("L").ColumnWidth = 2.5546875

$spacer = $ws.Range("G2:G10,L2:L10")
$spacer.Interior.Color = 14277081

# Row 2
$ws.Range("H2:I2").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C2)), 0), SUBSTITUTE(ANCHORARRAY(C2), "-", "") * -1,ANCHORARRAY(C2) + 0)'
$ws.Range("M2").Formula = '=IFERROR(SUM(ANCHORARRAY(H2)),0)'

# Row 3
$ws.Range("H3:J3").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C3)), 0), SUBSTITUTE(ANCHORARRAY(C3), "-", "") * -1,ANCHORARRAY(C3) + 0)'
$ws.Range("M3").Formula = '=IFERROR(SUM(ANCHORARRAY(H3)),0)'

# Row 4
$ws.Range("H4:K4").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C4)), 0), SUBSTITUTE(ANCHORARRAY(C4), "-", "") * -1,ANCHORARRAY(C4) + 0)'
$ws.Range("M4").Formula = '=IFERROR(SUM(ANCHORARRAY(H4)),0)'

# Row 5
$ws.Range("H5").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C5)), 0), SUBSTITUTE(ANCHORARRAY(C5), "-", "") * -1,ANCHORARRAY(C5) + 0)'
$ws.Range("M5").Formula = '=IFERROR(SUM(ANCHORARRAY(H5)),0)'

# Row 6
$ws.Range("H6").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C6)), 0), SUBSTITUTE(ANCHORARRAY(C6), "-", "") * -1,ANCHORARRAY(C6) + 0)'
$ws.Range("M6").Formula = '=IFERROR(SUM(ANCHORARRAY(H6)),0)'

# Row 7 (errors)
$ws.Range("H7").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C7)), 0), SUBSTITUTE(ANCHORARRAY(C7), "-", "") * -1,ANCHORARRAY(C7) + 0)'
$ws.Range("M7").Formula = '=IFERROR(SUM(ANCHORARRAY(H7)),0)'

# Row 8
$ws.Range("H8:K8").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C8)), 0), SUBSTITUTE(ANCHORARRAY(C8), "-", "") * -1,ANCHORARRAY(C8) + 0)'
$ws.Range("M8").Formula = '=IFERROR(SUM(ANCHORARRAY(H8)),0)'

# Row 9
$ws.Range("H9:K9").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C9)), 0), SUBSTITUTE(ANCHORARRAY(C9), "-", "") * -1,ANCHORARRAY(C9) + 0)'
$ws.Range("M9").Formula = '=IFERROR(SUM(ANCHORARRAY(H9)),0)'

# Row 10
$ws.Range("H10:K10").FormulaArray = '=IF(IFERROR(FIND("-",ANCHORARRAY(C10)), 0), SUBSTITUTE(ANCHORARRAY(C10), "-", "") * -1,ANCHORARRAY(C10) + 0)'
$ws.Range("M10").Formula = '=IFERROR(SUM(ANCHORARRAY(H10)),0)'

$ws.Range("D13").Select()

Write-Host "Done"
